$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cellRef, $val)
    $cell = $ws.Range($cellRef)
    if ($val -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = $origStyle
    } else {
        $cell.Value = $val
    }
}

$rows = @(
    @{ Row=2; D='26.840.08'; E='  -1.22%  ' },
    @{ Row=3; D='1.857.10'; E='  -0.76%  ' },
    @{ Row=4; E='  -0.15%  ' },
    @{ Row=5; D='304.84'; E='  -0.85%  ' },
    @{ Row=6; E='  -0.10%  ' },
    @{ Row=7; D='0.5073'; E='  -1.71%  ' },
    @{ Row=8; D='0.3654'; E='  -2.70%  ' },
    @{ Row=9; D='0.07182'; E='  +0.28%  ' },
    @{ Row=10; D='0.8896'; E='  +0.43%  ' },
    @{ Row=11; D='20.69'; E='  -0.80%  ' },
    @{ Row=12; D='0.07533'; E='  -0.61%  ' },
    @{ Row=13; D='1.853.38'; E='  -1.11%  ' },
    @{ Row=14; D='91.66'; E='  +2.44%  ' },
    @{ Row=15; D='5.231'; E='  -1.97%  ' },
    @{ Row=16; E='  -0.12%  ' },
    @{ Row=17; D='0.000008535'; E='  -0.19%  ' },
    @{ Row=18; D='14.06'; E='  -0.99%  ' },
    @{ Row=19; E='  -0.17%  ' },
    @{ Row=20; D='26.888.37'; E='  -1.20%  ' },
    @{ Row=21; D='5.017'; E='  -0.37%  ' },
    @{ Row=22; D='2.089.55'; E='  -1.40%  ' },
    @{ Row=23; D='10.32'; E='  -3.12%  ' },
    @{ Row=24; D='6.446'; E='  -0.61%  ' },
    @{ Row=25; E='  -3.62%  ' },
    @{ Row=26; D='1.803'; E='  -2.53%  ' },
    @{ Row=27; E='  -1.20%  ' },
    @{ Row=28; D='2.052'; E='  -5.73%  ' },
    @{ Row=29; D='112.99'; E='  -0.20%  ' },
    @{ Row=30; D='4.642'; E='  -2.50%  ' },
    @{ Row=31; E='  -0.78%  ' },
    @{ Row=32; D='0.09226'; E='  +2.12%  ' },
    @{ Row=33; E='  -1.50%  ' },
    @{ Row=34; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='3.073'; E='  -0.68%  ' },
    @{ Row=35; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.7347'; E='  -3.09%  ' },
    @{ Row=36; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.148'; E='  -2.77%  ' },
    @{ Row=37; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='3.215'; E='  +5.89%  ' },
    @{ Row=38; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.02007'; E='  -1.80%  ' },
    @{ Row=39; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='2.465'; E='  -1.89%  ' },
    @{ Row=40; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.074'; E='  -0.74%  ' },
    @{ Row=41; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.5308'; E='  -2.18%  ' },
    @{ Row=42; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='118.19'; E='  +2.68%  ' },
    @{ Row=43; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='6.500'; E='  -2.74%  ' },
    @{ Row=44; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='8.360'; E='  -2.43%  ' },
    @{ Row=45; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1472'; E='  -1.26%  ' },
    @{ Row=46; B='Decentraland'; C='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D='0.4635'; E='  -1.28%  ' },
    @{ Row=47; B='PaxDollar'; C='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D='0.9998'; E='  -0.10%  ' },
    @{ Row=48; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='9.956'; E='  -2.08%  ' },
    @{ Row=49; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='1.559'; E='  -1.26%  ' },
    @{ Row=50; B='Elrond'; C='https://coinranking.com/coin/omwkOTglq+elrond-egld'; D='36.94'; E='  +1.14%  ' },
    @{ Row=51; B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='62.87'; E='  -3.22%  ' }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    if ($r.ContainsKey("B")) { Set-TextCell "B$rowNum" $r.B }
    if ($r.ContainsKey("C")) { Set-TextCell "C$rowNum" $r.C }
    if ($r.ContainsKey("D")) { Set-TextCell "D$rowNum" $r.D }
    if ($r.ContainsKey("E")) { Set-TextCell "E$rowNum" $r.E }
}
